$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 46054
$ws.Range("B2").Value = 2.62
$ws.Range("C2").Value = 1.13
$ws.Range("D2").Value = 0.59
$ws.Range("E2").Value = 0.14
$ws.Range("F2").Value = 0.09
$ws.Range("G2").Value = 0.09
$ws.Range("H2").Value = 0.15
$ws.Range("I2").Value = 0.31
$ws.Range("J2").Value = 0.54
$ws.Range("K2").Value = 0.99
$ws.Range("L2").Value = 1.33
$ws.Range("M2").Value = 1.79
$ws.Range("N2").Value = 1.81
$ws.Range("O2").Value = 2.05
$ws.Range("P2").Value = 2.05
$ws.Range("Q2").Value = 1.75
$ws.Range("R2").Value = 1.55
$ws.Range("S2").Value = 1.1
$ws.Range("T2").Value = 0.41
$ws.Range("U2").Value = 1.15
$ws.Range("V2").Value = 1.3
$ws.Range("W2").Value = 1.23
$ws.Range("X2").Value = 0.72
$ws.Range("Y2").Value = 0.06
$ws.Range("Z2").Value = 1.04
$ws.Range("AA2").Value = "12h-16h"
$ws.Range("AB2").Value = 1.92
$ws.Range("AC2").Value = "12h-14h"
$ws.Range("AD2").Value = 1.93
$ws.Range("AE2").Value = "14h-16h"
$ws.Range("AF2").Value = 1.9
$ws.Range("AG2").Value = "2h-23h"
